$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the agenda text for the "Time to Code" block at row 7 (25 minute coding block)
# and the "Time for Retro" block at row 10, per the commit "Updated agenda and retro".
$ws.Range("B7").Value = "Time to Code with Red Green Refactor; do full cycles of TDD"
$ws.Range("B10").Value = "Time for Retro and Red Green Refactor Start Time Verification"

# Move the selection to B10 as in the saved workbook state.
$ws.Range("B10").Select() | Out-Null
